# Auto-generated Excel COM-interop script to apply odds updates
# Commit message: Atualizando o arquivo XLSX
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 4.2
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 17
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.95
$ws.Range("AI4").Value = 23
$ws.Range("AO4").Value = 10
$ws.Range("AT4").Value = 2.38
$ws.Range("AX4").Value = 6.5
$ws.Range("BA4").Value = 126
$ws.Range("G4").Value = 1.8
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.5
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("U4").Value = 2.25
$ws.Range("V4").Value = 1.57
$ws.Range("X4").Value = 7
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 13
$ws.Range("O6").Value = 1.13
$ws.Range("P6").Value = 5.2
$ws.Range("AA7").Value = 10
$ws.Range("AB7").Value = 23
$ws.Range("AC7").Value = 19
$ws.Range("AD7").Value = 11
$ws.Range("AE7").Value = 19
$ws.Range("AF7").Value = 51
$ws.Range("AG7").Value = 201
$ws.Range("AH7").Value = 26
$ws.Range("AL7").Value = 51
$ws.Range("AN7").Value = 3.5
$ws.Range("AP7").Value = 15
$ws.Range("AR7").Value = 34
$ws.Range("AS7").Value = 81
$ws.Range("AT7").Value = 4
$ws.Range("AU7").Value = 9
$ws.Range("AW7").Value = 351
$ws.Range("AZ7").Value = 34
$ws.Range("BA7").Value = 151
$ws.Range("BB7").Value = 126
$ws.Range("BC7").Value = 201
$ws.Range("G7").Value = 1.27
$ws.Range("I7").Value = 8
$ws.Range("K7").Value = 2.75
$ws.Range("L7").Value = 7.5
$ws.Range("M7").Value = 1.02
$ws.Range("N7").Value = 19
$ws.Range("O7").Value = 1.13
$ws.Range("P7").Value = 6
$ws.Range("Q7").Value = 1.44
$ws.Range("R7").Value = 2.7
$ws.Range("S7").Value = 1.22
$ws.Range("T7").Value = 4
$ws.Range("U7").Value = 1.73
$ws.Range("V7").Value = 2
$ws.Range("W7").Value = 10
$ws.Range("X7").Value = 8
$ws.Range("Y7").Value = 9.5
$ws.Range("Z7").Value = 9
$ws.Range("AH12").Value = 6.5
$ws.Range("AJ12").Value = 8.5
$ws.Range("AN12").Value = 9
$ws.Range("G12").Value = 10
$ws.Range("U12").Value = 2.05
$ws.Range("V12").Value = 1.7
$ws.Range("W12").Value = 21
$ws.Range("AA14").Value = 28
$ws.Range("AC14").Value = 9.5
$ws.Range("AD14").Value = 6.1
$ws.Range("AG14").Value = 350
$ws.Range("AH14").Value = 7.5
$ws.Range("AI14").Value = 10.5
$ws.Range("AJ14").Value = 8.5
$ws.Range("AK14").Value = 21
$ws.Range("AL14").Value = 17.5
$ws.Range("AM14").Value = 26
$ws.Range("AN14").Value = 5.2
$ws.Range("AO14").Value = 18
$ws.Range("AQ14").Value = 90
$ws.Range("AR14").Value = 120
$ws.Range("AT14").Value = 2.42
$ws.Range("AX14").Value = 4.05
$ws.Range("AY14").Value = 11.75
$ws.Range("AZ14").Value = 20
$ws.Range("BA14").Value = 50
$ws.Range("G14").Value = 3.35
$ws.Range("H14").Value = 3.1
$ws.Range("I14").Value = 2.15
$ws.Range("J14").Value = 3.8
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 2.82
$ws.Range("O14").Value = 1.28
$ws.Range("P14").Value = 3.05
$ws.Range("Q14").Value = 1.83
$ws.Range("R14").Value = 1.78
$ws.Range("T14").Value = 2.47
$ws.Range("W14").Value = 11
$ws.Range("X14").Value = 19.5
$ws.Range("Y14").Value = 11.25
$ws.Range("Z14").Value = 50
$ws.Range("AA15").Value = 16
$ws.Range("AB15").Value = 26
$ws.Range("AD15").Value = 6.5
$ws.Range("AE15").Value = 14
$ws.Range("AH15").Value = 11
$ws.Range("AI15").Value = 19.5
$ws.Range("AJ15").Value = 11.75
$ws.Range("AK15").Value = 50
$ws.Range("AL15").Value = 30
$ws.Range("AM15").Value = 35
$ws.Range("AN15").Value = 3.85
$ws.Range("AO15").Value = 10.25
$ws.Range("AP15").Value = 19
$ws.Range("AQ15").Value = 40
$ws.Range("AR15").Value = 75
$ws.Range("AT15").Value = 2.52
$ws.Range("AU15").Value = 7.1
$ws.Range("AX15").Value = 5.3
$ws.Range("AY15").Value = 19
$ws.Range("AZ15").Value = 25
$ws.Range("BA15").Value = 100
$ws.Range("G15").Value = 1.98
$ws.Range("H15").Value = 3.35
$ws.Range("I15").Value = 3.5
$ws.Range("J15").Value = 2.62
$ws.Range("K15").Value = 2.07
$ws.Range("L15").Value = 3.9
$ws.Range("N15").Value = 9.75
$ws.Range("O15").Value = 1.27
$ws.Range("P15").Value = 3.1
$ws.Range("Q15").Value = 1.82
$ws.Range("S15").Value = 1.39
$ws.Range("T15").Value = 2.55
$ws.Range("W15").Value = 7.4
$ws.Range("X15").Value = 9.5
$ws.Range("Y15").Value = 8.5
$ws.Range("Z15").Value = 18
$ws.Range("AC17").Value = 11
$ws.Range("AL17").Value = 51
$ws.Range("AM17").Value = 51
$ws.Range("AQ17").Value = 23
$ws.Range("AS17").Value = 151
$ws.Range("AU17").Value = 9
$ws.Range("BC17").Value = 301
$ws.Range("G17").Value = 1.48
$ws.Range("I17").Value = 6.5
$ws.Range("K17").Value = 2.3
$ws.Range("L17").Value = 6.5
$ws.Range("M17").Value = 1.05
$ws.Range("N17").Value = 11
$ws.Range("Q17").Value = 1.83
$ws.Range("R17").Value = 2.03
$ws.Range("AY18").Value = 67
$ws.Range("BA18").Value = 451
$ws.Range("X18").Value = 6
$ws.Range("AA19").Value = 29
$ws.Range("AD19").Value = 6
$ws.Range("AI19").Value = 9.5
$ws.Range("AS19").Value = 251
$ws.Range("G19").Value = 3.4
$ws.Range("I19").Value = 2.15
$ws.Range("O19").Value = 1.36
$ws.Range("P19").Value = 3
$ws.Range("Q19").Value = 2.2
$ws.Range("R19").Value = 1.65
$ws.Range("U19").Value = 1.95
$ws.Range("V19").Value = 1.8
$ws.Range("W19").Value = 8.5
$ws.Range("AE23").Value = 17
$ws.Range("AG23").Value = 151
$ws.Range("AH23").Value = 10
$ws.Range("AX23").Value = 3.6
$ws.Range("AY23").Value = 6.5
$ws.Range("BA23").Value = 17
$ws.Range("G23").Value = 7.5
$ws.Range("H23").Value = 5
$ws.Range("I23").Value = 1.36
$ws.Range("J23").Value = 6.5
$ws.Range("Y23").Value = 21
